# Prototype is finished
# Update quiz answers on the lakersquiz sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lakersquiz")

# Row 8: "Why did Alex Caruso shave his head (Hint: Awareness)?"
$ws.Range("D8").Value = "Donation to his former school; Donation to Church for Tax Deductible;Donation for Democratic Party"

# Row 23: "Who works as a Executive Security Agent for Lakers (As of Sep 2021)?"
$ws.Range("D23").Value = "Frank Thompson;Roland Dore;Andrew Brook"

# Row 15: "Was Game 6 against Sacramento Kings during 2002 Playoffs rigged?"
$ws.Range("D15").Value = "No;No;No"

# Row 7: "What date did Alex Caruso sign with Oklahoma City Thunder?"
$ws.Range("D7").Value = "8/16/2015;11/13/2016;1/13/2017"

# Row 3: "What's Lebron's Nickname?"
$ws.Range("C3").Value = "The King"
$ws.Range("D3").Value = "The Queen;The Qing;Black Mamba"

# Update the last selected cell to reflect the final cursor position after editing.
$ws.Range("D4").Select()
